$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 349.16666
$ws.Range("J9").Value = 556.6667
$ws.Range("L9").Value = 556.6667
$ws.Range("N9").Value = -894.6667
$ws.Range("H41").Value = 436.81818
$ws.Range("I41").Value = 449.83334
$ws.Range("J41").Value = 421.2
$ws.Range("K41").Value = 449.83334
$ws.Range("L41").Value = 421.2
$ws.Range("M41").Value = -9.833340000000021
$ws.Range("N41").Value = -1301.2
$ws.Range("H53").Value = 436.48
$ws.Range("I53").Value = 551.8125
$ws.Range("J53").Value = 231.44444
$ws.Range("K53").Value = 551.8125
$ws.Range("L53").Value = 231.44444
$ws.Range("M53").Value = 85.1875
$ws.Range("N53").Value = -1505.44444
$ws.Range("H76").Value = 2282.1428
$ws.Range("I76").Value = 1995
$ws.Range("K76").Value = 1995
$ws.Range("M76").Value = -1680
$ws.Range("H79").Value = 2282.1428
$ws.Range("I79").Value = 1995
$ws.Range("K79").Value = 1995
$ws.Range("M79").Value = -903
$ws.Range("H86").Value = 44790.668
$ws.Range("I86").Value = 13186.5
$ws.Range("K86").Value = 13186.5
$ws.Range("M86").Value = -12063.5
$ws.Range("H89").Value = 44790.668
$ws.Range("I89").Value = 13186.5
$ws.Range("K89").Value = 65932.5
$ws.Range("M89").Value = -60316.5
$ws.Range("H138").Value = 2020.5769
$ws.Range("J138").Value = 3324.348
$ws.Range("L138").Value = 9973.044
$ws.Range("N138").Value = -20253.044

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1832.375
$ws.Range("I61").Value = 1833.2069
$ws.Range("K61").Value = 1833.2069
$ws.Range("M61").Value = -1621.2069
$ws.Range("H93").Value = 34965.5
$ws.Range("J93").Value = 34965.5
$ws.Range("L93").Value = 34965.5
$ws.Range("N93").Value = -39957.5
$ws.Range("H122").Value = 2773.077
$ws.Range("I122").Value = 2550.0908
$ws.Range("K122").Value = 7650.2724
$ws.Range("M122").Value = -5200.2724
$ws.Range("H132").Value = 1513.7142
$ws.Range("I132").Value = 1368.683
$ws.Range("J132").Value = 2257
$ws.Range("K132").Value = 4106.049
$ws.Range("L132").Value = 6771
$ws.Range("M132").Value = -1576.049
$ws.Range("N132").Value = -11831
$ws.Range("H136").Value = 1832.375
$ws.Range("I136").Value = 1833.2069
$ws.Range("K136").Value = 5499.620699999999
$ws.Range("M136").Value = -2949.620699999999
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 16599.666
$ws.Range("J49").Value = 16599.666
$ws.Range("L49").Value = 16599.666
$ws.Range("N49").Value = -17077.666
$ws.Range("H94").Value = 559.3077
$ws.Range("I94").Value = 512.5
$ws.Range("K94").Value = 512.5
$ws.Range("M94").Value = -61.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2134.2632
$ws.Range("I31").Value = 1912.4375
$ws.Range("J31").Value = 3317.3333
$ws.Range("K31").Value = 1912.4375
$ws.Range("L31").Value = 3317.3333
$ws.Range("M31").Value = -1617.4375
$ws.Range("N31").Value = -3907.3333
$ws.Range("H34").Value = 2134.2632
$ws.Range("I34").Value = 1912.4375
$ws.Range("J34").Value = 3317.3333
$ws.Range("K34").Value = 1912.4375
$ws.Range("L34").Value = 3317.3333
$ws.Range("M34").Value = -1710.4375
$ws.Range("N34").Value = -3721.3333
$ws.Range("H74").Value = 46666.668
$ws.Range("J74").Value = 46666.668
$ws.Range("L74").Value = 46666.668
$ws.Range("N74").Value = -48414.668
$ws.Range("H77").Value = 46666.668
$ws.Range("J77").Value = 46666.668
$ws.Range("L77").Value = 140000.004
$ws.Range("N77").Value = -148736.004
$ws.Range("H122").Value = 3547.946
$ws.Range("I122").Value = 2055.9285
$ws.Range("K122").Value = 6167.7855
$ws.Range("M122").Value = -3717.7855
$ws.Range("H134").Value = 9746.286
$ws.Range("I134").Value = 10182.846
$ws.Range("J134").Value = 9036.875
$ws.Range("K134").Value = 30548.538
$ws.Range("L134").Value = 27110.625
$ws.Range("M134").Value = -28013.538
$ws.Range("N134").Value = -32180.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 672.46155
$ws.Range("I23").Value = 755.8
$ws.Range("J23").Value = 620.375
$ws.Range("K23").Value = 2267.4
$ws.Range("L23").Value = 1861.125
$ws.Range("M23").Value = -2032.4
$ws.Range("N23").Value = -2331.125
$ws.Range("H121").Value = 12821804
$ws.Range("J121").Value = 2007.25
$ws.Range("L121").Value = 6021.75
$ws.Range("N121").Value = -8641.75
$ws.Range("H140").Value = 2033
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 111111170
$ws.Range("I2").Value = 194444510
$ws.Range("J2").Value = 26.333334
$ws.Range("K2").Value = 194444510
$ws.Range("L2").Value = 26.333334
$ws.Range("M2").Value = -194444397
$ws.Range("N2").Value = -252.333334
$ws.Range("H3").Value = 6250730
$ws.Range("I3").Value = 1667306.4
$ws.Range("K3").Value = 1667306.4
$ws.Range("M3").Value = -1667190.4
$ws.Range("H9").Value = 13250.75
$ws.Range("I9").Value = 4331.6665
$ws.Range("K9").Value = 4331.6665
$ws.Range("M9").Value = -4161.6665
$ws.Range("H70").Value = 17233.572
$ws.Range("I70").Value = 28408.75
$ws.Range("K70").Value = 28408.75
$ws.Range("M70").Value = -28138.75
$ws.Range("H73").Value = 17233.572
$ws.Range("I73").Value = 28408.75
$ws.Range("K73").Value = 28408.75
$ws.Range("M73").Value = -27472.75
$ws.Range("H122").Value = 772532.3
$ws.Range("I122").Value = 1432430.2
$ws.Range("K122").Value = 4297290.6
$ws.Range("M122").Value = -4294840.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 595.2857
$ws.Range("I9").Value = 283.6
$ws.Range("J9").Value = 1374.5
$ws.Range("K9").Value = 283.6
$ws.Range("L9").Value = 1374.5
$ws.Range("M9").Value = -59.60000000000002
$ws.Range("N9").Value = -1822.5
$ws.Range("H16").Value = 4202.6
$ws.Range("I16").Value = 3565.5
$ws.Range("J16").Value = 6751
$ws.Range("K16").Value = 3565.5
$ws.Range("L16").Value = 6751
$ws.Range("M16").Value = -3395.5
$ws.Range("N16").Value = -7091
$ws.Range("H55").Value = 824.9666999999999
$ws.Range("I55").Value = 298.16666
$ws.Range("K55").Value = 298.16666
$ws.Range("M55").Value = -125.16666
$ws.Range("H138").Value = 87746.75
$ws.Range("J138").Value = 87746.75
$ws.Range("L138").Value = 87746.75
$ws.Range("N138").Value = -98026.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9769.625
$ws.Range("J41").Value = 9011.666999999999
$ws.Range("L41").Value = 9011.666999999999
$ws.Range("N41").Value = -9791.666999999999
$ws.Range("H113").Value = 1057.2858
$ws.Range("J113").Value = 2003
$ws.Range("L113").Value = 6009
$ws.Range("N113").Value = -10349
